# Sửa lại file thông điệp như ban đầu (mã 3 và 11):
# In the "Server phát bài cho 4 Client" bullet, the run "TenClient%" that
# sits between "3%" and "La1%La2%....%La13" is dropped, so the message
# reverts to "3%La1%La2%....%La13".
$d = $word.ActiveDocument

$d.Content.Find.Execute("TenClient%", $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 2) | Out-Null
